# Auto-generated edit script applying the diff to Garuda_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1951
$ws.Range("I6").Value = 51
$ws.Range("J6").Value = 2901
$ws.Range("K6").Value = 153
$ws.Range("L6").Value = 8703
$ws.Range("M6").Value = -41
$ws.Range("N6").Value = -8927
$ws.Range("H137").Value = 1154.079
$ws.Range("I137").Value = 934
$ws.Range("J137").Value = 1398.6111
$ws.Range("K137").Value = 2802
$ws.Range("L137").Value = 4195.8333
$ws.Range("M137").Value = -252
$ws.Range("N137").Value = -9295.8333
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14765.69
$ws.Range("I32").Value = 14598.639
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 14598.639
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -14311.639
$ws.Range("N32").Value = -20574
$ws.Range("H61").Value = 2016.6666
$ws.Range("I61").Value = 1735.7333
$ws.Range("J61").Value = 2484.889
$ws.Range("K61").Value = 1735.7333
$ws.Range("L61").Value = 2484.889
$ws.Range("M61").Value = -1523.7333
$ws.Range("N61").Value = -2908.889
$ws.Range("H110").Value = 2094.36
$ws.Range("I110").Value = 2115.0435
$ws.Range("J110").Value = 1856.5
$ws.Range("K110").Value = 2115.0435
$ws.Range("L110").Value = 1856.5
$ws.Range("M110").Value = -70.04350000000022
$ws.Range("N110").Value = -5946.5
$ws.Range("H119").Value = 20000
$ws.Range("I119").Value = 10000
$ws.Range("J119").Value = 25000
$ws.Range("K119").Value = 10000
$ws.Range("L119").Value = 25000
$ws.Range("N119").Value = -34676
$ws.Range("M119").Value = -5162
$ws.Range("H124").Value = 14500
$ws.Range("J124").Value = 14500
$ws.Range("L124").Value = 14500
$ws.Range("N124").Value = -24320
$ws.Range("H136").Value = 2016.6666
$ws.Range("I136").Value = 1735.7333
$ws.Range("J136").Value = 2484.889
$ws.Range("K136").Value = 5207.199900000001
$ws.Range("L136").Value = 7454.667
$ws.Range("M136").Value = -2657.199900000001
$ws.Range("N136").Value = -12554.667
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3062.3914
$ws.Range("I105").Value = 2033.5385
$ws.Range("J105").Value = 4399.9
$ws.Range("K105").Value = 2033.5385
$ws.Range("L105").Value = 4399.9
$ws.Range("M105").Value = -286.5385000000001
$ws.Range("N105").Value = -7893.9
$ws.Range("H132").Value = 88333.336
$ws.Range("J132").Value = 88333.336
$ws.Range("L132").Value = 88333.336
$ws.Range("N132").Value = -98453.336
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1056.625
$ws.Range("I16").Value = 1047.5
$ws.Range("J16").Value = 1065.75
$ws.Range("K16").Value = 1047.5
$ws.Range("L16").Value = 1065.75
$ws.Range("M16").Value = -760.5
$ws.Range("N16").Value = -1639.75
$ws.Range("H100").Value = 39900
$ws.Range("J100").Value = 39900
$ws.Range("L100").Value = 39900
$ws.Range("N100").Value = -42064
$ws.Range("H105").Value = 779.7143
$ws.Range("I105").Value = 571.8
$ws.Range("J105").Value = 1299.5
$ws.Range("K105").Value = 571.8
$ws.Range("L105").Value = 1299.5
$ws.Range("M105").Value = 1175.2
$ws.Range("N105").Value = -4793.5
$ws.Range("H113").Value = 1056.625
$ws.Range("I113").Value = 1047.5
$ws.Range("J113").Value = 1065.75
$ws.Range("K113").Value = 1047.5
$ws.Range("L113").Value = 1065.75
$ws.Range("M113").Value = 1122.5
$ws.Range("N113").Value = -5405.75
$ws.Range("H141").Value = 47228.344
$ws.Range("I141").Value = 18314.666
$ws.Range("J141").Value = 54771.043
$ws.Range("K141").Value = 18314.666
$ws.Range("L141").Value = 54771.043
$ws.Range("M141").Value = -13134.666
$ws.Range("N141").Value = -65131.043
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H80").Value = 3675
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 3675
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H92").Value = 642.8570999999999
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 683.3333
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 2049.9999
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -4545.9999
$ws.Range("H131").Value = 631.86
$ws.Range("I131").Value = 253.4
$ws.Range("J131").Value = 835.6462
$ws.Range("K131").Value = 760.2
$ws.Range("L131").Value = 2506.9386
$ws.Range("M131").Value = 4279.8
$ws.Range("N131").Value = -12586.9386
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 21999
$ws.Range("J131").Value = 21999
$ws.Range("L131").Value = 21999
$ws.Range("N131").Value = -32079
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 54989
$ws.Range("J127").Value = 54989
$ws.Range("L127").Value = 54989
$ws.Range("N127").Value = -64909
$ws.Range("H132").Value = 6524.3486
$ws.Range("I132").Value = 8541.25
$ws.Range("K132").Value = 25623.75
$ws.Range("M132").Value = -23093.75
$ws.Range("H133").Value = 26756.908
$ws.Range("J133").Value = 26756.908
$ws.Range("L133").Value = 26756.908
$ws.Range("N133").Value = -31816.908
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29450
$ws.Range("J70").Value = 29450
$ws.Range("L70").Value = 29450
$ws.Range("N70").Value = -30080
$ws.Range("H73").Value = 29450
$ws.Range("J73").Value = 29450
$ws.Range("L73").Value = 29450
$ws.Range("N73").Value = -31634
$ws.Range("H126").Value = 31255928
$ws.Range("I126").Value = 52639292
$ws.Range("J126").Value = 3320.1538
$ws.Range("K126").Value = 157917876
$ws.Range("L126").Value = 9960.4614
$ws.Range("M126").Value = -157915406
$ws.Range("N126").Value = -14900.4614
$ws.Range("H131").Value = 37905
$ws.Range("J131").Value = 37905
$ws.Range("L131").Value = 37905
